# Error Calculations and Plots
# Update missing-data pattern: fill/clear a few "D" (column E) values,
# and drop two rows (RM 232, SC 92) from the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in previously-missing values
$ws.Range("E3").Value = -5.7
$ws.Range("E21").Value = -8.699999999999999

# Clear values that are now treated as missing
$ws.Range("E5").Value = ""
$ws.Range("E23").Value = ""

# Remove the "RM 232" row entirely (was row 26)
$ws.Rows.Item(26).Delete()

# Remove the "SC 92" row entirely (now shifted up to row 27)
$ws.Rows.Item(27).Delete()

# Fill in the value for the row that is now "SC 193" (row 32 after the shifts)
$ws.Range("E32").Value = -6.4
